$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "dugwell"
$ws.Range("C1").Value = "spring"
$ws.Range("D1").Value = "piped_home"
$ws.Range("E1").Value = "piped_yard_tap"
$ws.Range("F1").Value = "cheap_commercial"
$ws.Range("G1").Value = "surface_water"
$ws.Range("H1").Value = "thirthy_min_less_travel"
